$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.718.33"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$ws.Range("D3").Value = "2.289.22"
$ws.Range("E3").Value = "  -1.90%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.25"
$ws.Range("E5").Value = "  +1.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.18"
$ws.Range("E6").Value = "  -0.18%  "

# Row 7
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -1.75%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.66"
$ws.Range("E10").Value = "  +1.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.95"
$ws.Range("E12").Value = "  -1.62%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +1.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.65"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("D15").Value = "2.633.76"
$ws.Range("E15").Value = "  -1.33%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -0.46%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.289.19"
$ws.Range("E17").Value = "  -1.37%  "

# Row 18
$ws.Range("D18").Value = "43.665.38"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000112"
$ws.Range("E19").Value = "  +3.52%  "

# Row 20
$ws.Range("E20").Value = "  -2.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.14"
$ws.Range("E21").Value = "  +0.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").Value = "  +9.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.73"
$ws.Range("E23").Value = "  -2.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  -5.97%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("E25").Value = "  +7.21%  "

# Row 26
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.32"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.45"
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.70"
$ws.Range("E29").Value = "  +0.28%  "

# Row 30
$ws.Range("E30").Value = "  -2.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.74"
$ws.Range("E31").Value = "  +1.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.19"
$ws.Range("E32").Value = "  -2.70%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("E34").Value = "  -0.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.64"
$ws.Range("E36").Value = "  +6.29%  "

# Row 37
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0352"
$ws.Range("E38").Value = "  -1.39%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  +5.06%  "

# Row 40
$ws.Range("E40").Value = "  +1.17%  "

# Row 41
$ws.Range("E41").Value = "  -2.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.35"
$ws.Range("E42").Value = "  +2.03%  "

# Row 43
$ws.Range("E43").Value = "  -2.38%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.56"
$ws.Range("E44").Value = "  +4.63%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.73"
$ws.Range("E45").Value = "  -4.27%  "

# Row 46
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.20"
$ws.Range("E46").Value = "  -3.06%  "

# Row 48
$ws.Range("E48").Value = "  -0.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.48"
$ws.Range("E49").Value = "  -3.12%  "

# Row 50
$ws.Range("E50").Value = "  +12.21%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.432"
$ws.Range("E51").Value = "  +4.25%  "
